# MAI_holdings.xlsx - "Add files via upload"
#
# Update the "as of" date in the confidential disclosure footer
# (2021-05-27 -> 2021-05-28) and refresh the Weight / Percent Change
# figures for the six model-holding rows (D2:E7).
#
# The worksheet is protected (legacy password hash "D382"), so cells are
# locked by default. Unlock only the cells being touched, write the new
# values, then re-lock them so the sheet's protected state is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Disclosure footer text (shared string) -------------------------------
$ws.Range("A10").Locked = $false
$ws.Range("A10").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-28 for illustrative purposes only and are subject to change."
$ws.Range("A10").Locked = $true
# Editing the multi-line text triggers an auto row-height bump; restore the
# row to its original (default-height, non-custom) state.
$ws.Rows(10).AutoFit()

# --- Weight (D) / Percent Change (E) for rows 2-7 --------------------------
$rng = $ws.Range("D2:E7")
$rng.Locked = $false

$ws.Range("D2").Value2 = 0.4778972637973691
$ws.Range("E2").Value2 = 0.004674717569146836

$ws.Range("D3").Value2 = 0.3396565661514342
$ws.Range("E3").Value2 = 0.001406337896118437

$ws.Range("D4").Value2 = 0.09731033975944933
$ws.Range("E4").Value2 = 0.00845814977973558

$ws.Range("D5").Value2 = 0.05340739289126227
$ws.Range("E5").Value2 = -0.0003439972480218101

$ws.Range("D6").Value2 = 0.03172843740048502
$ws.Range("E6").Value2 = -0.00579038795599307

$ws.Range("D7").Value2 = 1
$ws.Range("E7").Value2 = 0.003332680106802322

$rng.Locked = $true
